$wb = $excel.ActiveWorkbook

# Rename sheets for excel/csv consistency
$wb.Worksheets.Item("field_mapping").Name = "fields"
$wb.Worksheets.Item("value_mapping").Name = "values"

$ws = $wb.Worksheets.Item("data")
$ws.Range("D2:D28").FormulaR1C1 = "=RC[-1]-RC[-2]"

$wb.Worksheets.Item("values").Activate()

foreach ($sh in $wb.Worksheets) {
    Write-Host $sh.Name
}
